$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells hold text values (inline strings) in the canonical XML,
# so force Text number format before assignment to prevent Excel from
# auto-converting numeric-looking strings (e.g. "1.002", "26.10") into numbers.


# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.962.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.40%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.70%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.60%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4036"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.82%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4073"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.002"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.67%  "

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.73"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.31%  "

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.473"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08822"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.10"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.492"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.52%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.052"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001353"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.62%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.775.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.93%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.59"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.99%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07155"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.245"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.83%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.942.41"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.36%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.907"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.493"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +24.73%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.87"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.34%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "145.61"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.84%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.256"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.949.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.76%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.240"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +13.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08856"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.36%  "

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03208"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.36%  "

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.385"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.024"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2845"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.63%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8428"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.48%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.87"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09369"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.46%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.13"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.33%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.474"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.722"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7448"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.07%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.247"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.393"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.89%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.51"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.66%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08373"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.55%  "
